# Add sentiment values (column E) for tweet rows 886-1001 (session 1 annotation batch)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sentimentData = @(
    @{Row=886; Value="NEUTRAL"},
    @{Row=887; Value="NEUTRAL"},
    @{Row=888; Value="NEUTRAL"},
    @{Row=889; Value="NEUTRAL"},
    @{Row=890; Value="NEUTRAL"},
    @{Row=891; Value="NEGATIVE"},
    @{Row=892; Value="NEUTRAL"},
    @{Row=893; Value="NEUTRAL"},
    @{Row=894; Value="NEUTRAL"},
    @{Row=895; Value="NEGATIVE"},
    @{Row=896; Value="MIXED"},
    @{Row=897; Value="NEGATIVE"},
    @{Row=898; Value="NEGATIVE"},
    @{Row=899; Value="NEGATIVE"},
    @{Row=900; Value="NEGATIVE"},
    @{Row=901; Value="NEGATIVE"},
    @{Row=902; Value="NEUTRAL"},
    @{Row=903; Value="NEUTRAL"},
    @{Row=904; Value="NEGATIVE"},
    @{Row=905; Value="NEUTRAL"},
    @{Row=906; Value="POSITIVE"},
    @{Row=907; Value="NEGATIVE"},
    @{Row=908; Value="NEGATIVE"},
    @{Row=909; Value="NEGATIVE"},
    @{Row=910; Value="NEUTRAL"},
    @{Row=911; Value="NEGATIVE"},
    @{Row=912; Value="NEUTRAL"},
    @{Row=913; Value="POSITIVE"},
    @{Row=914; Value="NEGATIVE"},
    @{Row=915; Value="NEGATIVE"},
    @{Row=916; Value="NEUTRAL"},
    @{Row=917; Value="NEUTRAL"},
    @{Row=918; Value="NEUTRAL"},
    @{Row=919; Value="MIXED"},
    @{Row=920; Value="NEGATIVE"},
    @{Row=921; Value="NEGATIVE"},
    @{Row=922; Value="NEUTRAL"},
    @{Row=923; Value="NEGATIVE"},
    @{Row=924; Value="NEGATIVE"},
    @{Row=925; Value="NEGATIVE"},
    @{Row=926; Value="NEGATIVE"},
    @{Row=927; Value="NEUTRAL"},
    @{Row=928; Value="NEUTRAL"},
    @{Row=929; Value="NEUTRAL"},
    @{Row=930; Value="NEGATIVE"},
    @{Row=931; Value="NEGATIVE"},
    @{Row=932; Value="NEUTRAL"},
    @{Row=933; Value="NEUTRAL"},
    @{Row=934; Value="NEGATIVE"},
    @{Row=935; Value="NEGATIVE"},
    @{Row=936; Value="NEGATIVE"},
    @{Row=937; Value="NEUTRAL"},
    @{Row=938; Value="NEGATIVE"},
    @{Row=939; Value="NEUTRAL"},
    @{Row=940; Value="MIXED"},
    @{Row=941; Value="NEUTRAL"},
    @{Row=942; Value="NEUTRAL"},
    @{Row=943; Value="NEUTRAL"},
    @{Row=944; Value="NEGATIVE"},
    @{Row=945; Value="NEGATIVE"},
    @{Row=946; Value="NEGATIVE"},
    @{Row=947; Value="MIXED"},
    @{Row=948; Value="NEGATIVE"},
    @{Row=949; Value="NEGATIVE"},
    @{Row=950; Value="NEGATIVE"},
    @{Row=951; Value="NEUTRAL"},
    @{Row=952; Value="NEUTRAL"},
    @{Row=953; Value="NEGATIVE"},
    @{Row=954; Value="MIXED"},
    @{Row=955; Value="NEUTRAL"},
    @{Row=956; Value="NEGATIVE"},
    @{Row=957; Value="POSITIVE"},
    @{Row=958; Value="POSITIVE"},
    @{Row=959; Value="NEGATIVE"},
    @{Row=960; Value="NEGATIVE"},
    @{Row=961; Value="NEUTRAL"},
    @{Row=962; Value="NEGATIVE"},
    @{Row=963; Value="NEUTRAL"},
    @{Row=964; Value="POSITIVE"},
    @{Row=965; Value="NEUTRAL"},
    @{Row=966; Value="NEGATIVE"},
    @{Row=967; Value="NEUTRAL"},
    @{Row=968; Value="NEGATIVE"},
    @{Row=969; Value="NEUTRAL"},
    @{Row=970; Value="POSITIVE"},
    @{Row=971; Value="POSITIVE"},
    @{Row=972; Value="NEGATIVE"},
    @{Row=973; Value="MIXED"},
    @{Row=974; Value="NEUTRAL"},
    @{Row=975; Value="NEUTRAL"},
    @{Row=976; Value="POSITIVE"},
    @{Row=977; Value="NEGATIVE"},
    @{Row=978; Value="POSITIVE"},
    @{Row=979; Value="NEGATIVE"},
    @{Row=980; Value="NEGATIVE"},
    @{Row=981; Value="NEGATIVE"},
    @{Row=982; Value="NEUTRAL"},
    @{Row=983; Value="NEUTRAL"},
    @{Row=984; Value="NEUTRAL"},
    @{Row=985; Value="NEUTRAL"},
    @{Row=986; Value="NEUTRAL"},
    @{Row=987; Value="NEUTRAL"},
    @{Row=988; Value="NEUTRAL"},
    @{Row=989; Value="POSITIVE"},
    @{Row=990; Value="NEGATIVE"},
    @{Row=991; Value="NEGATIVE"},
    @{Row=992; Value="POSITIVE"},
    @{Row=993; Value="NEUTRAL"},
    @{Row=994; Value="POSITIVE"},
    @{Row=995; Value="NEUTRAL"},
    @{Row=996; Value="NEGATIVE"},
    @{Row=997; Value="NEGATIVE"},
    @{Row=998; Value="NEGATIVE"},
    @{Row=999; Value="NEGATIVE"},
    @{Row=1000; Value="NEGATIVE"},
    @{Row=1001; Value="NEGATIVE"}
)

foreach ($item in $sentimentData) {
    $ws.Cells.Item($item.Row, 5).Value = $item.Value
}

# Update the sheet view to reflect the new scroll position / active selection
$ws.Application.ActiveWindow.ScrollRow = 989
$ws.Range("E1001").Select()
